$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.805.87"
$ws.Range("E2").Value = "'  +1.15%  "

$ws.Range("D3").Value = "'3.403.91"
$ws.Range("E3").Value = "'  +0.29%  "

$ws.Range("E4").Value = "'  -0.16%  "

$ws.Range("D5").Value = "'560.96"

$ws.Range("D6").Value = "'176.47"
$ws.Range("E6").Value = "'  +1.14%  "

$ws.Range("E7").Value = "'  +0.92%  "

$ws.Range("D8").Value = "'3.396.25"
$ws.Range("E8").Value = "'  +0.27%  "

$ws.Range("E9").Value = "'  +0.02%  "

$ws.Range("D10").Value = "'0.176"
$ws.Range("E10").Value = "'  +6.09%  "

$ws.Range("E11").Value = "'  +0.58%  "

$ws.Range("D12").Value = "'53.70"
$ws.Range("E12").Value = "'  -1.64%  "

$ws.Range("D13").Value = "'0.0000280"
$ws.Range("E13").Value = "'  +1.04%  "

$ws.Range("D14").Value = "'9.23"
$ws.Range("E14").Value = "'  +1.15%  "

$ws.Range("D15").Value = "'3.940.08"
$ws.Range("E15").Value = "'  -0.13%  "

$ws.Range("D16").Value = "'18.32"
$ws.Range("E16").Value = "'  +0.13%  "

$ws.Range("D17").Value = "'3.410.58"
$ws.Range("E17").Value = "'  +0.27%  "

$ws.Range("E18").Value = "'  +0.94%  "

$ws.Range("D19").Value = "'65.770.89"
$ws.Range("E19").Value = "'  +0.96%  "

$ws.Range("D20").Value = "'11.88"
$ws.Range("E20").Value = "'  -0.25%  "

$ws.Range("E21").Value = "'  +0.57%  "

$ws.Range("D22").Value = "'480.93"
$ws.Range("E22").Value = "'  +2.09%  "

$ws.Range("D23").Value = "'4.94"
$ws.Range("E23").Value = "'  -1.28%  "

$ws.Range("B24").Value = "'InternetComputer(DFINITY)"
$ws.Range("C24").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'14.36"
$ws.Range("E24").Value = "'  +5.31%  "

$ws.Range("B25").Value = "'PancakeSwap"
$ws.Range("C25").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'4.12"
$ws.Range("E25").Value = "'  -0.34%  "

$ws.Range("D26").Value = "'89.38"
$ws.Range("E26").Value = "'  +3.60%  "

$ws.Range("D27").Value = "'2.92"
$ws.Range("E27").Value = "'  +1.88%  "

$ws.Range("D28").Value = "'10.68"
$ws.Range("E28").Value = "'  -1.54%  "

$ws.Range("D29").Value = "'8.75"
$ws.Range("E29").Value = "'  -0.80%  "

$ws.Range("E30").Value = "'  +2.50%  "

$ws.Range("D31").Value = "'6.59"
$ws.Range("E31").Value = "'  -1.56%  "

$ws.Range("D32").Value = "'63.73"
$ws.Range("E32").Value = "'  +6.00%  "

$ws.Range("E33").Value = "'  -0.34%  "

$ws.Range("D34").Value = "'571.53"
$ws.Range("E34").Value = "'  -2.30%  "

$ws.Range("E35").Value = "'  -0.47%  "

$ws.Range("E36").Value = "'  -0.03%  "

$ws.Range("D37").Value = "'3.66"
$ws.Range("E37").Value = "'  +3.15%  "

$ws.Range("E38").Value = "'  +1.05%  "

$ws.Range("D39").Value = "'35.91"
$ws.Range("E39").Value = "'  -0.07%  "

$ws.Range("E40").Value = "'  +0.38%  "

$ws.Range("D41").Value = "'0.0₃0746"
$ws.Range("E41").Value = "'  -0.79%  "

$ws.Range("D42").Value = "'3.091.11"
$ws.Range("E42").Value = "'  -0.42%  "

$ws.Range("D43").Value = "'2.81"
$ws.Range("E43").Value = "'  -1.72%  "

$ws.Range("E44").Value = "'  +0.54%  "

$ws.Range("E45").Value = "'  +0.47%  "

$ws.Range("D46").Value = "'2.46"
$ws.Range("E46").Value = "'  -2.41%  "

$ws.Range("D47").Value = "'3.17"
$ws.Range("E47").Value = "'  -1.76%  "

$ws.Range("D48").Value = "'0.998"
$ws.Range("E48").Value = "'  -0.17%  "

$ws.Range("D49").Value = "'140.18"
$ws.Range("E49").Value = "'  +2.33%  "

$ws.Range("E50").Value = "'  -0.14%  "

$ws.Range("D51").Value = "'8.47"
$ws.Range("E51").Value = "'  +1.63%  "
